$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Save" in H1, matching the formatting of the other headers (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data values in column H for rows 2 and 3
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
